# "TP 7/ejercicio2.xlsx" — add a "k =" label and the mean of the M column
# (the random-number sample average) below the random-numbers table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 53: label in L53, average formula in M53.
$ws.Range("L53").Value = "k ="
$ws.Range("M53").Formula = "=AVERAGE(M2:M51)"

# Move the active selection to where the user ended up after the edit.
$null = $ws.Range("S8").Select()
